# limit_up_history.xlsx / "晋级率" sheet: append the 2025-05-16 promotion-rate
# breakdown as a new block of 7 rows (541-547), continuing the existing
# alternating white/gray row-block shading (previous block, 2025-05-15, was
# white, so this block uses the gray F2F2F2 fill).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 541
$lastRow = 547

# Shade the whole new block gray (RGB F2F2F2 -> decimal 15921906) to match the
# alternating block-fill pattern already used throughout the sheet.
$ws.Range("A${firstRow}:P${lastRow}").Interior.Color = 15921906

# Columns A, D-J hold text in the source sheet (dates, "x进y" labels and
# percentages formatted as strings) rather than real dates/numbers, so force
# a text number format before writing them to avoid Excel auto-converting
# "2025-05-16" to a date serial or "26.92%" to a numeric percentage.
$textCols = @("A", "D", "E", "F", "G", "H", "I", "J")
foreach ($col in $textCols) {
    $ws.Range("${col}${firstRow}:${col}${lastRow}").NumberFormat = "@"
}

# Row 541
$ws.Range("A541").Value = "2025-05-16"
$ws.Range("D541").Value = "1进2"
$ws.Range("E541").Value = "26.92%"
$ws.Range("F541").Value = "61.54%"
$ws.Range("G541").Value = "38.46%"
$ws.Range("H541").Value = "38.89%"
$ws.Range("I541").Value = "83.33%"
$ws.Range("J541").Value = "16.67%"
$ws.Range("B541").Value = 26
$ws.Range("C541").Value = 18
$ws.Range("K541").Value = 7
$ws.Range("L541").Value = 16
$ws.Range("M541").Value = 10
$ws.Range("N541").Value = 7
$ws.Range("O541").Value = 15
$ws.Range("P541").Value = 3

# Row 542
$ws.Range("A542").Value = "2025-05-16"
$ws.Range("D542").Value = "2进3"
$ws.Range("E542").Value = "11.11%"
$ws.Range("F542").Value = "55.56%"
$ws.Range("G542").Value = "44.44%"
$ws.Range("H542").Value = "0.00%"
$ws.Range("I542").Value = "100.00%"
$ws.Range("J542").Value = "0.00%"
$ws.Range("B542").Value = 9
$ws.Range("C542").Value = 3
$ws.Range("K542").Value = 1
$ws.Range("L542").Value = 5
$ws.Range("M542").Value = 4
$ws.Range("N542").Value = 0
$ws.Range("O542").Value = 3
$ws.Range("P542").Value = 0

# Row 543
$ws.Range("A543").Value = "2025-05-16"
$ws.Range("D543").Value = "3进4"
$ws.Range("E543").Value = "71.43%"
$ws.Range("F543").Value = "71.43%"
$ws.Range("G543").Value = "28.57%"
$ws.Range("H543").Value = "100.00%"
$ws.Range("I543").Value = "100.00%"
$ws.Range("J543").Value = "0.00%"
$ws.Range("B543").Value = 7
$ws.Range("C543").Value = 5
$ws.Range("K543").Value = 5
$ws.Range("L543").Value = 5
$ws.Range("M543").Value = 2
$ws.Range("N543").Value = 5
$ws.Range("O543").Value = 5
$ws.Range("P543").Value = 0

# Row 544
$ws.Range("A544").Value = "2025-05-16"
$ws.Range("D544").Value = "4进5"
$ws.Range("E544").Value = "0.00%"
$ws.Range("F544").Value = "0.00%"
$ws.Range("G544").Value = "100.00%"
$ws.Range("H544").Value = "0.00%"
$ws.Range("I544").Value = "0.00%"
$ws.Range("J544").Value = "100.00%"
$ws.Range("B544").Value = 1
$ws.Range("C544").Value = 1
$ws.Range("K544").Value = 0
$ws.Range("L544").Value = 0
$ws.Range("M544").Value = 1
$ws.Range("N544").Value = 0
$ws.Range("O544").Value = 0
$ws.Range("P544").Value = 1

# Row 545
$ws.Range("A545").Value = "2025-05-16"
$ws.Range("D545").Value = "5进6"
$ws.Range("E545").Value = "0.00%"
$ws.Range("F545").Value = "0.00%"
$ws.Range("G545").Value = "100.00%"
$ws.Range("H545").Value = "0.00%"
$ws.Range("I545").Value = "0.00%"
$ws.Range("J545").Value = "0.00%"
$ws.Range("B545").Value = 2
$ws.Range("C545").Value = 0
$ws.Range("K545").Value = 0
$ws.Range("L545").Value = 0
$ws.Range("M545").Value = 2
$ws.Range("N545").Value = 0
$ws.Range("O545").Value = 0
$ws.Range("P545").Value = 0

# Row 546
$ws.Range("A546").Value = "2025-05-16"
$ws.Range("D546").Value = "7进8"
$ws.Range("E546").Value = "100.00%"
$ws.Range("F546").Value = "100.00%"
$ws.Range("G546").Value = "0.00%"
$ws.Range("H546").Value = "100.00%"
$ws.Range("I546").Value = "100.00%"
$ws.Range("J546").Value = "0.00%"
$ws.Range("B546").Value = 1
$ws.Range("C546").Value = 1
$ws.Range("K546").Value = 1
$ws.Range("L546").Value = 1
$ws.Range("M546").Value = 0
$ws.Range("N546").Value = 1
$ws.Range("O546").Value = 1
$ws.Range("P546").Value = 0

# Row 547
$ws.Range("A547").Value = "2025-05-16"
$ws.Range("D547").Value = "8进9"
$ws.Range("E547").Value = "100.00%"
$ws.Range("F547").Value = "100.00%"
$ws.Range("G547").Value = "0.00%"
$ws.Range("H547").Value = "100.00%"
$ws.Range("I547").Value = "100.00%"
$ws.Range("J547").Value = "0.00%"
$ws.Range("B547").Value = 1
$ws.Range("C547").Value = 1
$ws.Range("K547").Value = 1
$ws.Range("L547").Value = 1
$ws.Range("M547").Value = 0
$ws.Range("N547").Value = 1
$ws.Range("O547").Value = 1
$ws.Range("P547").Value = 0
